$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (B11): Right marks multiplier 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row (B12): total right-answer score 39 -> 65
$ws.Range("B12").Value = 65

# Update the Corr/Total marks text (E12): "35/84" -> "65/140"
$ws.Range("E12").Value = "65/140"
